# Restructured code for multiple trial options:
# adds three new per-trial sheets (TAG9, TAG10, TAG12) and appends their
# aggregated rows to the "summary" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Append the new per-trial detail sheets at the end of the workbook.
# ---------------------------------------------------------------------

# TAG9 (-> xl/worksheets/sheet9.xml)
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$tag9 = $wb.Worksheets.Add($null, $last)
$tag9.Name = "TAG9"
$tag9.Range("A1").Value = "Trial"
$tag9.Range("B1").Value = "Reaction Time"
$tag9.Range("C1").Value = "Accuracy"

$tag9Rows = @(
    ,@("('folder_3', 'ball_3')", "3.28s", "correct")
    ,@("('atest', 'ball_1 - Copy')", "1.94s", "correct")
    ,@("('folder_4', 'ball_4')", "1.95s", "correct")
    ,@("('atest', 'ball_1')", "1.90s", "correct")
    ,@("('folder_2', 'ball_2')", "1.96s", "correct")
)
$r = 2
foreach ($row in $tag9Rows) {
    $tag9.Cells.Item($r, 1).Value = $row[0]
    $tag9.Cells.Item($r, 2).Value = $row[1]
    $tag9.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# TAG10 (-> xl/worksheets/sheet10.xml)
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$tag10 = $wb.Worksheets.Add($null, $last)
$tag10.Name = "TAG10"
$tag10.Range("A1").Value = "Trial"
$tag10.Range("B1").Value = "Reaction Time"
$tag10.Range("C1").Value = "Accuracy"

$tag10Rows = @(
    ,@("('folder_3', 'ball_3')", "1.96s", "correct")
    ,@("('atest', 'ball_1')", "0.63s", "correct")
    ,@("('folder_4', 'ball_4')", "0.43s", "correct")
    ,@("('atest', 'ball_1 - Copy')", "0.40s", "correct")
    ,@("('folder_2', 'ball_2')", "0.46s", "correct")
)
$r = 2
foreach ($row in $tag10Rows) {
    $tag10.Cells.Item($r, 1).Value = $row[0]
    $tag10.Cells.Item($r, 2).Value = $row[1]
    $tag10.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# TAG12 (-> xl/worksheets/sheet11.xml)
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$tag12 = $wb.Worksheets.Add($null, $last)
$tag12.Name = "TAG12"
$tag12.Range("A1").Value = "Trial"
$tag12.Range("B1").Value = "Reaction Time"
$tag12.Range("C1").Value = "Accuracy"

$tag12Rows = @(
    ,@("('folder_3', 'ball_3')", "2.03s", "correct")
    ,@("('folder_4', 'ball_4')", "0.82s", "incorrect")
    ,@("('folder_2', 'ball_2')", "0.00s", "incorrect")
    ,@("('atest', 'ball_1 - Copy')", "0.00s", "incorrect")
    ,@("('atest', 'ball_1')", "0.56s", "correct")
)
$r = 2
foreach ($row in $tag12Rows) {
    $tag12.Cells.Item($r, 1).Value = $row[0]
    $tag12.Cells.Item($r, 2).Value = $row[1]
    $tag12.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Append the corresponding aggregate rows to the "summary" sheet.
# ---------------------------------------------------------------------

$summary = $wb.Worksheets.Item("summary")

$summaryRows = @(
    ,@("TAG9",  "2.21s", "5/5", "2.21s", "5/5")
    ,@("TAG10", "0.78s", "5/5", "0.78s", "5/5")
    ,@("TAG",   "0.00s", "0/0", "0.00s", "0/0")
    ,@("TAG12", "0.68s", "2/5", "0.68s", "2/5")
)

$r = 10
foreach ($row in $summaryRows) {
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
    $summary.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
